$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting the existing "INDIA - ISL" (row 3)
# and "ROMANIA - LIGA 1" (row 4) rows down to rows 4 and 5 respectively.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new GREECE - SUPER LEAGUE match.
$ws.Range("A3").Value = "YHOF5rBm"
$ws.Range("B3").Value = "25/11/2024"
$ws.Range("C3").Value = "13:00"
$ws.Range("D3").Value = "GREECE - SUPER LEAGUE"
$ws.Range("E3").Value = "Atromitos"
$ws.Range("F3").Value = "Levadiakos"
$ws.Range("G3").Value = 1.83
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.75
$ws.Range("J3").Value = 2.6
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 5.5
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 7.5
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 15
$ws.Range("AA3").Value = 19
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 6.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 501
$ws.Range("AH3").Value = 10
$ws.Range("AI3").Value = 23
$ws.Range("AJ3").Value = 17
$ws.Range("AK3").Value = 51
$ws.Range("AL3").Value = 41
$ws.Range("AM3").Value = 51
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 10
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.38
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 6
$ws.Range("AX3").Value = 29
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 101
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 401
# Apply odds updates to the ROMANIA - LIGA 1 match, now on row 5.
$ws.Range("G5").Value = 2.35
$ws.Range("H5").Value = 2.88
$ws.Range("I5").Value = 3.2
$ws.Range("J5").Value = 3.2
$ws.Range("W5").Value = 6.5
$ws.Range("AA5").Value = 21
$ws.Range("AB5").Value = 34
$ws.Range("AE5").Value = 17
$ws.Range("AH5").Value = 8.5
$ws.Range("AJ5").Value = 13
$ws.Range("AW5").Value = 5
$ws.Range("AY5").Value = 34
$ws.Range("AZ5").Value = 67
$ws.Range("BB5").Value = 301
# Append the new UKRAINE - PREMIER LEAGUE match as row 6.
$ws.Range("A6").Value = "zRQv9vQQ"
$ws.Range("B6").Value = "25/11/2024"
$ws.Range("C6").Value = "13:00"
$ws.Range("D6").Value = "UKRAINE - PREMIER LEAGUE"
$ws.Range("E6").Value = "Polissya Zhytomyr"
$ws.Range("F6").Value = "FK Zorya Luhansk"
$ws.Range("G6").Value = 1.7
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 4.9
$ws.Range("J6").Value = 2.25
$ws.Range("K6").Value = 2.07
$ws.Range("L6").Value = 5.2
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 9.050000000000001
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 2.67
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.62
$ws.Range("S6").Value = 1.42
$ws.Range("T6").Value = 2.45
$ws.Range("U6").Value = 1.98
$ws.Range("V6").Value = 1.65
$ws.Range("W6").Value = 5.7
$ws.Range("X6").Value = 7
$ws.Range("Y6").Value = 8.25
$ws.Range("Z6").Value = 12.5
$ws.Range("AA6").Value = 15
$ws.Range("AB6").Value = 32
$ws.Range("AC6").Value = 8.25
$ws.Range("AD6").Value = 6.7
$ws.Range("AE6").Value = 18.5
$ws.Range("AF6").Value = 110
$ws.Range("AG6").Value = 900
$ws.Range("AH6").Value = 11.5
$ws.Range("AI6").Value = 28
$ws.Range("AJ6").Value = 16.5
$ws.Range("AK6").Value = 100
$ws.Range("AL6").Value = 60
$ws.Range("AM6").Value = 65
$ws.Range("AN6").Value = 3.35
$ws.Range("AO6").Value = 8.25
$ws.Range("AP6").Value = 19.5
$ws.Range("AQ6").Value = 28
$ws.Range("AR6").Value = 65
$ws.Range("AS6").Value = 300
$ws.Range("AT6").Value = 2.42
$ws.Range("AU6").Value = 7.9
$ws.Range("AV6").Value = 80
$ws.Range("AW6").Value = 6.5
$ws.Range("AX6").Value = 30
$ws.Range("AY6").Value = 37
$ws.Range("AZ6").Value = 200
$ws.Range("BA6").Value = 250
$ws.Range("BB6").Value = 500
$ws.Range("BC6").Value = 81
$ws.Range("BD6").Value = 81